# Applies the cryptos.xlsx price/volume/coin-order refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.257.28"
$ws.Range("E2").Value = "  -2.97%  "
$ws.Range("D3").Value = "3.849.30"
$ws.Range("E3").Value = "  -3.04%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'601.20"
$ws.Range("E5").Value = "  -1.98%  "
$ws.Range("E6").Value = "  -1.53%  "
$ws.Range("D7").Value = "3.847.99"
$ws.Range("E7").Value = "  -3.16%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  -1.86%  "
$ws.Range("E10").Value = "  -5.30%  "
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("E12").Value = "  -3.03%  "
$ws.Range("E13").Value = "  +0.76%  "
$ws.Range("D14").Value = "'36.98"
$ws.Range("E14").Value = "  -3.64%  "
$ws.Range("D15").Value = "4.494.15"
$ws.Range("D16").Value = "3.847.63"
$ws.Range("E16").Value = "  -2.76%  "
$ws.Range("D17").Value = "68.299.99"
$ws.Range("E17").Value = "  -2.79%  "
$ws.Range("D18").Value = "'18.08"
$ws.Range("E18").Value = "  +0.89%  "
$ws.Range("D19").Value = "'7.35"
$ws.Range("E19").Value = "  -4.21%  "
$ws.Range("E20").Value = "  -0.71%  "
$ws.Range("D21").Value = "'10.76"
$ws.Range("E21").Value = "  -3.41%  "
$ws.Range("D22").Value = "'466.03"
$ws.Range("E22").Value = "  -7.42%  "
$ws.Range("D23").Value = "'0.732"
$ws.Range("E23").Value = "  -1.62%  "
$ws.Range("E24").Value = "  -6.01%  "
$ws.Range("E25").Value = "  -3.45%  "
$ws.Range("E26").Value = "  -3.52%  "
$ws.Range("D27").Value = "'12.04"
$ws.Range("E27").Value = "  -3.74%  "
$ws.Range("D28").Value = "'10.01"
$ws.Range("E28").Value = "  -2.90%  "
$ws.Range("E30").Value = "  -1.45%  "
$ws.Range("D31").Value = "3.998.42"
$ws.Range("E31").Value = "  -2.94%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'2.31"
$ws.Range("E32").Value = "  -5.52%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "'7.59"
$ws.Range("E33").Value = "  -4.52%  "
$ws.Range("D34").Value = "'31.23"
$ws.Range("E34").Value = "  -4.10%  "
$ws.Range("E35").Value = "  -1.92%  "
$ws.Range("D36").Value = "3.812.75"
$ws.Range("E36").Value = "  -3.11%  "
$ws.Range("E37").Value = "  -4.01%  "
$ws.Range("E38").Value = "  +10.10%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.141"
$ws.Range("E39").Value = "  -0.33%  "
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D40").Value = "'1.02"
$ws.Range("E40").Value = "  -2.17%  "
$ws.Range("E41").Value = "  -4.78%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").Value = "'0.314"
$ws.Range("E44").Value = "  -5.09%  "
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").Value = "'420.14"
$ws.Range("E45").Value = "  -4.90%  "
$ws.Range("B46").Value = "FLOKI"
$ws.Range("C46").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D46").Value = "'0.000294"
$ws.Range("E46").Value = "  +5.80%  "
$ws.Range("D48").Value = "'8.60"
$ws.Range("E48").Value = "  -0.89%  "
$ws.Range("D49").Value = "'47.07"
$ws.Range("E49").Value = "  -2.65%  "
$ws.Range("D50").Value = "'143.04"
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("D51").Value = "'26.14"
$ws.Range("E51").Value = "  -0.03%  "
